$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cyclic rotation of rows 16, 17, 18 (row16 <- old row18, row17 <- old row16, row18 <- old row17)
# with updated "Taxonsorteringsordning" (column B) values.

$ws.Range("A16").Value = 112178654
$ws.Range("B16").Value = 89834
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 658
$ws.Range("F16").Value = "Rosenticka"
$ws.Range("G16").Value = "Rhodofomes roseus"
$ws.Range("H16").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q16").Value = 618387
$ws.Range("R16").Value = 6904851

$ws.Range("A17").Value = 112178652
$ws.Range("B17").Value = 90826
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 4366
$ws.Range("F17").Value = "Skarp dropptaggsvamp"
$ws.Range("G17").Value = "Hydnellum peckii"
$ws.Range("H17").Value = "Banker"
$ws.Range("Q17").Value = 618476
$ws.Range("R17").Value = 6905002

$ws.Range("A18").Value = 112178651
$ws.Range("B18").Value = 86371
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 4412
$ws.Range("F18").Value = "Äggvaxskivling"
$ws.Range("G18").Value = "Hygrophorus karstenii"
$ws.Range("H18").Value = "Sacc. & Cub."
$ws.Range("Q18").Value = 618388
$ws.Range("R18").Value = 6904949
